# Swap the species-identification data between row 2 and row 4.
# Columns involved: A (Id), B (Taxonsorteringsordning), E (TaxonId),
# F (Artnamn), G (Vetenskapligt namn), H (Auktor).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H")

foreach ($col in $cols) {
    $addr2 = "${col}2"
    $addr4 = "${col}4"

    $cell2 = $ws.Range($addr2)
    $cell4 = $ws.Range($addr4)

    $val2 = $cell2.Value2
    $val4 = $cell4.Value2

    $cell2.Value2 = $val4
    $cell4.Value2 = $val2
}
